$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with two new columns P and Q, matching the style of the
# existing header cells (B1:O1 use style index 1, i.e. the same style as A1 etc.)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Fill data rows 2-25 for the two new columns with 0, no special style (same as existing
# data cells in columns B:O for those rows).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}
